$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.011.89"
$ws.Range("E2").Value = "  +0.75%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.310.90"
$ws.Range("E3").Value = "  +0.45%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.10"
$ws.Range("E5").Value = "  -1.56%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.79"
$ws.Range("E6").Value = "  +2.55%  "

# Row 7
$ws.Range("E7").Value = "  -0.31%  "

# Row 8
$ws.Range("E8").Value = "  +0.21%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.605"
$ws.Range("E9").Value = "  +0.03%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.07"
$ws.Range("E10").Value = "  +1.09%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0913"
$ws.Range("E11").Value = "  +0.64%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.36"
$ws.Range("E12").Value = "  -2.22%  "

# Row 13
$ws.Range("E13").Value = "  -0.32%  "

# Row 14
$ws.Range("E14").Value = "  -1.17%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.30"
$ws.Range("E15").Value = "  -0.43%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.657.85"
$ws.Range("E16").Value = "  +0.53%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.311.44"
$ws.Range("E17").Value = "  +0.96%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.866.31"
$ws.Range("E18").Value = "  +0.65%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.48"
$ws.Range("E19").Value = "  -0.99%  "

# Row 20
$ws.Range("E20").Value = "  -0.85%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.37"
$ws.Range("E21").Value = "  -3.00%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.55"
$ws.Range("E22").Value = "  -0.66%  "

# Row 23
$ws.Range("E23").Value = "  -1.11%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.06"
$ws.Range("E24").Value = "  -0.37%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.24"
$ws.Range("E25").Value = "  -0.23%  "

# Row 26
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  -0.04%  "

# Row 27
$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.76"
$ws.Range("E27").Value = "  +17.15%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.96"
$ws.Range("E28").Value = "  +0.31%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.50"
$ws.Range("E30").Value = "  +3.47%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.30"
$ws.Range("E31").Value = "  -1.49%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.59"
$ws.Range("E32").Value = "  -0.01%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0869"
$ws.Range("E33").Value = "  -1.59%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.78"
$ws.Range("E34").Value = "  +8.81%  "

# Row 35
$ws.Range("E35").Value = "  -0.68%  "

# Row 36
$ws.Range("E36").Value = "  +1.87%  "

# Row 37
$ws.Range("E37").Value = "  -1.24%  "

# Row 38
$ws.Range("E38").Value = "  +1.14%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.81"
$ws.Range("E39").Value = "  +3.60%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.62"
$ws.Range("E40").Value = "  -3.56%  "

# Row 41
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.58"
$ws.Range("E41").Value = "  -0.95%  "

# Row 42
$ws.Range("B42").Value = "BitcoinSV"
$ws.Range("C42").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "103.63"
$ws.Range("E42").Value = "  +7.88%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.79"
$ws.Range("E43").Value = "  +0.13%  "

# Row 44
$ws.Range("E44").Value = "  +0.72%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.23%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.27"
$ws.Range("E46").Value = "  -1.60%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "111.64"
$ws.Range("E47").Value = "  -4.76%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.692.76"
$ws.Range("E48").Value = "  +2.08%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.25"
$ws.Range("E49").Value = "  -4.88%  "

# Row 50
$ws.Range("E50").Value = "  +0.30%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.21"
$ws.Range("E51").Value = "  -1.48%  "
